$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G4").Value = "2016-09-06 11:00:45"

$wsZhCn.Range("H4").Value = "2016-09-06 11:00:35"
$wsZhCn.Range("K4").Value = "2016-09-06 11:01:26"

$wsDeDe.Range("K4").Value = "2016-09-06 11:01:34"
